$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.931.38'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  -2.29%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.753.76'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  -4.90%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.22'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  -9.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9994'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  -0.18%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5031'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  -5.63%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.74'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  -7.10%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2636'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').Value = '  -13.87%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06178'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  -10.39%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.750.13'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').Value = '  -5.74%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06942'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  -11.00%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.37'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  -16.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.475'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  -10.52%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.33'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  -13.96%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.5877'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  -21.60%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9988'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').Value = '  -0.25%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9989'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  -0.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.957.30'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  -2.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.67'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').Value = '  -16.73%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006760'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  -14.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.973.38'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  -5.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.063'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  -12.22%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.079'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  -13.44%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.103'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  -14.96%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '138.05'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  -3.71%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.538'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  -8.96%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.841'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  -16.54%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.88'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  -12.52%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.61'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  -7.53%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.768'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  -11.93%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08124'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  -7.83%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.471'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  -14.85%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04474'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  -7.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9980'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  -0.27%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.620'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  -10.60%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9949'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  -12.47%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6037'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  -17.29%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.682'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').Value = '  -13.70%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.937'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  -16.12%  '

$ws.Range('B41').Value = 'VeChain'

$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01535'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  -10.90%  '

$ws.Range('B42').Value = 'Quant'

$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '103.90'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  -4.35%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9991'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  -0.20%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3814'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  -20.50%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.144'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  -12.59%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7331'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  -19.53%  '

$ws.Range('B47').Value = 'Cronos'

$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05257'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  -9.40%  '

$ws.Range('B48').Value = 'Algorand'

$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1109'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  -10.84%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.949'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  -20.82%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.22'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').Value = '  -13.44%  '

$ws.Range('E51').Value = '  -13.43%  '
